$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.807.30'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').Value = '3.406.82'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '407.62'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '128.26'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.80%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.631'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +6.42%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.729'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +7.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.141'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +15.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.40'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.140'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000213'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +61.08%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '3.959.89'
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.88'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +5.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.84'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +4.60%  '
$ws.Range('D17').Value = '3.417.50'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.05'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +9.13%  '
$ws.Range('E19').Value = '  +3.50%  '
$ws.Range('D20').Value = '61.772.99'
$ws.Range('E20').Value = '  -0.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '407.55'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +29.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '89.28'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +5.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.17'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.04'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.23'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +3.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '32.74'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +10.00%  '
$ws.Range('E27').Value = '  +0.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.56'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.60'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.73'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.118'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.81%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.171'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.88%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.82'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '43.10'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.21%  '
$ws.Range('E35').Value = '  +0.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0494'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '53.96'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +3.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.35'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.46%  '
$ws.Range('E40').Value = '  +6.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.91'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.60%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.309'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +5.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '140.15'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.96'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.04'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.06%  '
$ws.Range('E46').Value = '  +8.54%  '
$ws.Range('E47').Value = '  -1.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '21.83'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.17%  '
$ws.Range('D49').Value = '2.113.94'
$ws.Range('E49').Value = '  -0.82%  '
$ws.Range('E50').Value = '  +4.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.132'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +16.65%  '
